# "implementation grading sheet hw6 claimed"
#
# The grader is entering the "claimed" (actual) hw6 points in column C for
# the items where the full "possible" score already recorded in column B
# was earned: row 6 (accuracy of degree requirements) and rows 10-15
# (functionality items 5.2 import degree requirements .. 5.7 add student).
# Copying the graded B-column cell down into the matching C-column cell
# reproduces both the value and the cell's number format/style exactly
# like a grader pasting their score claim over from the possible-points
# column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$claimedRows = @(6, 10, 11, 12, 13, 14, 15)
foreach ($r in $claimedRows) {
    $src = $ws.Range("B$r")
    $dst = $ws.Range("C$r")
    $src.Copy($dst)
}

# Leave the cursor where the grader finished typing.
[void]$ws.Range("M15").Select()
